$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.797000000000001
$ws.Range("B9").Value = 6.484999999999999
$ws.Range("B18").Value = 6.351
$ws.Range("B20").Value = 6.667999999999999
